$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jänner")

# Log a first day's worth of work: 1.5h on CCS
$ws.Range("A4").Value = 42644
$ws.Range("B4").Value = 1.5
$ws.Range("C4").Value = "CCS zum Laufen bringen"

# A few days later: 4h, same topic + gyro sensor
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = "CCS zum Laufen bringen; Gyro Sensor auslesen"

# 1.5h researching the sensor board
$ws.Range("B8").Value = 1.5
$ws.Range("C8").Value = "Recherche Sensor Board"

[void]$ws.Range("B10").Select()
